# Commit: added new line for description of component table for a better view
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the wrapped description text
$ws.Columns.Item(2).ColumnWidth = 37.25

# Rewrite each description (column B) with an inserted line break, enable wrap text,
# and grow the row height to fit the now multi-line text.
$ws.Range("B2").Value = "A simple regression model that fits a linear`n relationship between the input features and the`n target variable."
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 58.0

$ws.Range("B3").Value = "A linear regression model with L2 regularization `nto prevent overfitting by penalizing large coefficients."
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 58.0

$ws.Range("B4").Value = "A linear regression model with L1 regularization `nthat encourages sparsity by setting some coefficients to zero."
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 58.0

$ws.Range("B5").Value = "A combination of L1 (Lasso) and L2 (Ridge) `nregularization to balance feature selection and shrinkage."
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 43.5

$ws.Range("B6").Value = "A linear regression model trained using`n Stochastic Gradient Descent (SGD) for large-scale datasets."
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 43.5

$ws.Range("B7").Value = "A linear classifier that estimates probabilities `nusing the logistic function for binary`nclassification."
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 43.5

$ws.Range("B8").Value = "A Ridge-regularized version of logistic `nregression to improve generalization."
$ws.Range("B8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 29.0

$ws.Range("B9").Value = "A linear classifier trained using Stochastic `nGradient Descent (SGD) for large-scale datasets."
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 43.5

$ws.Range("B10").Value = "A Support Vector Regression model with a `nradial basis function (RBF) kernel, capturing complex non-linear relationships."
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 43.5

$ws.Range("B11").Value = "A Support Vector Regression model using a `nlinear kernel for simple relationships."
$ws.Range("B11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 29.0

$ws.Range("B12").Value = "A Support Vector Regression model with a `npolynomial kernel to capture non-linear dependencies."
$ws.Range("B12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 43.5

$ws.Range("B13").Value = "A Support Vector Regression model with a `nsigmoid kernel, often used in specialized cases."
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 43.5

$ws.Range("B14").Value = "A Support Vector Classification model using`n a linear kernel for simple classification tasks."
$ws.Range("B14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 29.0

$ws.Range("B15").Value = "A Support Vector Classification model with `nan RBF kernel, capturing non-linear relationships."
$ws.Range("B15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 43.5

$ws.Range("B16").Value = "A Support Vector Classification model using `na polynomial kernel for complex decision boundaries."
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 43.5

$ws.Range("B17").Value = "A Support Vector Classification model with a `nsigmoid kernel, suitable for specific tasks."
$ws.Range("B17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 29.0

$ws.Range("B18").Value = "A non-linear regression model that splits data `ninto hierarchical decision rules."
$ws.Range("B18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 29.0

$ws.Range("B19").Value = "A non-linear classifier that splits data based on `nfeature values to make decisions."
$ws.Range("B19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 43.5

$ws.Range("B20").Value = "An ensemble method that combines multiple `ndecision trees to improve prediction accuracy."
$ws.Range("B20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 43.5

$ws.Range("B21").Value = "An ensemble method that aggregates multiple `ndecision trees to improve accuracy and reduce overfitting."
$ws.Range("B21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 58.0

$ws.Range("B22").Value = "A boosting model that sequentially improves `nweak decision trees for enhanced predictive performance."
$ws.Range("B22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 43.5

$ws.Range("B23").Value = "A boosting model that sequentially improves `nweak classifiers for better performance."
$ws.Range("B23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 29.0

$ws.Range("B24").Value = "A boosting method that adjusts weak models `niteratively to minimize errors."
$ws.Range("B24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 29.0

$ws.Range("B25").Value = "A boosting method that assigns more weight to `nmisclassified instances in successive models."
$ws.Range("B25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 58.0

$ws.Range("B26").Value = "An ensemble method that trains multiple `ninstances of a regression model on different data samples for robustness."
$ws.Range("B26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 43.5

$ws.Range("B27").Value = "An ensemble method that improves model `nstability by training classifiers on different random data subsets."
$ws.Range("B27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 43.5

$ws.Range("B28").Value = "A probabilistic classifier that assumes normal `ndistribution of feature values and independence between them."
$ws.Range("B28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 43.5

$ws.Range("B29").Value = "A variant of Naïve Bayes suitable for discrete `nfeatures like word counts in text classification."
$ws.Range("B29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 43.5

$ws.Range("B30").Value = "A Naïve Bayes model designed for binary/`nboolean feature data."
$ws.Range("B30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 29.0

$ws.Range("B31").Value = "A non-parametric regression model that `npredicts a target value based on the average of the k-nearest data points."
$ws.Range("B31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 43.5

$ws.Range("B32").Value = "A non-parametric classification model that `nassigns a label based on the majority vote of k-nearest data points."
$ws.Range("B32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 43.5

$ws.Range("B33").Value = "Standardizes features by removing the mean `nand scaling to unit variance (Z-score normalization)."
$ws.Range("B33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 43.5

$ws.Range("B34").Value = "Scales features to a fixed range (default [0,1]) `nby transforming each value proportionally."
$ws.Range("B34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 29.0

$ws.Range("B35").Value = "Scales each feature by its maximum absolute `nvalue, preserving the sign and keeping values between -1 and 1."
$ws.Range("B35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 43.5

$ws.Range("B36").Value = "Uses median and interquartile range (IQR) for `nscaling, making it resistant to outliers."
$ws.Range("B36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 29.0

$ws.Range("B37").Value = "Scales individual samples to have unit norm, `nuseful for distance-based models like KNN and SVM."
$ws.Range("B37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 43.5

$ws.Range("B38").Value = "Converts categorical labels into numeric `nvalues (e.g., ‘red’, ‘blue’, ‘green’ → 0, 1, 2)."
$ws.Range("B38").WrapText = $true
$ws.Rows.Item(38).RowHeight = 29.0

$ws.Range("B39").Value = "Converts categorical values into binary `ncolumns, creating a separate column for each unique category."
$ws.Range("B39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 43.5

$ws.Range("B40").Value = "Encodes categorical values as ordinal numbers `n(useful when categories have an inherent order)."
$ws.Range("B40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 58.0

$ws.Range("B41").Value = "Converts multi-class labels into binary (one-`nvs-rest) format for classification tasks."
$ws.Range("B41").WrapText = $true
$ws.Rows.Item(41).RowHeight = 29.0

$ws.Range("B42").Value = "Replaces missing values with a specified `nstrategy (e.g., mean, median, most frequent)."
$ws.Range("B42").WrapText = $true
$ws.Rows.Item(42).RowHeight = 29.0

$ws.Range("B43").Value = "Uses k-nearest neighbors to fill in missing `nvalues based on similar samples."
$ws.Range("B43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 29.0

$ws.Range("B44").Value = "Converts numeric values into binary format `nbased on a threshold (e.g., all values >0.5 become 1, else 0)."
$ws.Range("B44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 43.5

$ws.Range("B45").Value = "A general utility that trains a machine learning `nmodel on a dataset."
$ws.Range("B45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 29.0

$ws.Range("B46").Value = "Uses a trained model to make predictions on `nnew data."
$ws.Range("B46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 29.0

$ws.Range("B47").Value = "Fits a preprocessing step (e.g., scaling, `nencoding) to the dataset before applying transformations."
$ws.Range("B47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 43.5

$ws.Range("B48").Value = "Applies transformations like scaling, encoding, `nor feature engineering to modify the dataset."
$ws.Range("B48").WrapText = $true
$ws.Rows.Item(48).RowHeight = 43.5

$ws.Range("B49").Value = "A combined step that both fits and transforms `ndata in one go, useful in pipelines."
$ws.Range("B49").WrapText = $true
$ws.Rows.Item(49).RowHeight = 29.0

$ws.Range("B50").Value = "Divides data into separate subsets, such as `ntraining and testing sets."
$ws.Range("B50").WrapText = $true
$ws.Rows.Item(50).RowHeight = 29.0

$ws.Range("B51").Value = "Loads data from various sources (CSV, `ndatabase, API, etc.) for use in machine learning models."
$ws.Range("B51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 43.5

$ws.Range("B52").Value = "Saves processed data or model states into a `nstructured format (e.g., JSON, database, file system)."
$ws.Range("B52").WrapText = $true
$ws.Rows.Item(52).RowHeight = 43.5

$ws.Range("B53").Value = "Loads saved node data for reuse in machine `nlearning workflows."
$ws.Range("B53").WrapText = $true
$ws.Rows.Item(53).RowHeight = 29.0

$ws.Range("B54").Value = "Merges datasets or data frames based on `ncommon attributes, useful for data preprocessing."
$ws.Range("B54").WrapText = $true
$ws.Rows.Item(54).RowHeight = 43.5

$ws.Range("B55").Value = "Assesses model performance using metrics `nsuch as accuracy, precision, recall, RMSE, etc."
$ws.Range("B55").WrapText = $true
$ws.Rows.Item(55).RowHeight = 43.5

$ws.Range("B56").Value = "Splits data into training and testing sets to `nvalidate model performance."
$ws.Range("B56").WrapText = $true
$ws.Rows.Item(56).RowHeight = 29.0

# Restore view/selection state similar to the source edit
$ws.Range("A47").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B28").Select()

